# Apply the recorded edits to MT_10.xlsx:
#  - Sheet2 gains a small summary table (fairness stats per SCHED_Type)
#  - Sheet2 becomes the active/selected sheet (tabSelected + workbook activeTab)
#  - Sheet2's selection becomes H8:H9 (active cell H8)
#  - MT_10 (sheet1) loses tabSelected, and its selection grows from N56 to N56:P56
#  - Sheet2 column A gets a custom width

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("MT_10")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Populate Sheet2 with the fairness summary table ---
$ws2.Range("B1").Value = "BATCH"
$ws2.Range("C1").Value = "DEADLINE"
$ws2.Range("D1").Value = "FIFO"
$ws2.Range("E1").Value = "OTHER"
$ws2.Range("F1").Value = "RR"

$ws2.Range("A2").Value = "Fairness(QWT)"
$ws2.Range("B2").Value = 0.011
$ws2.Range("C2").Value = 0.29942611776530126
$ws2.Range("D2").Value = 7.617008861751441
$ws2.Range("E2").Value = 0.014
$ws2.Range("F2").Value = 0.2915904662364668

$ws2.Range("A3").Value = "Fairness(RT)"
$ws2.Range("B3").Value = 0.0790189850605536
$ws2.Range("C3").Value = 0.347690954728477
$ws2.Range("D3").Value = 0.039799497484264874
$ws2.Range("E3").Value = 0.05063595560468899
$ws2.Range("F3").Value = 0.12727922061357855

$ws2.Range("A4").Value = "Fairness(ET)"
$ws2.Range("B4").Value = 0.08075270893288987
$ws2.Range("C4").Value = 0.3900371777151506
$ws2.Range("D4").Value = 7.610758175109751
$ws2.Range("E4").Value = 0.04512205669071412
$ws2.Range("F4").Value = 0.1657860066471231

# Column A on Sheet2 is widened to fit the "Fairness(...)" labels
# (target stored width 16.375 chars; engine quantizes to 1/7ths, so this is
# the closest achievable setting).
$ws2.Columns.Item(1).ColumnWidth = 15.714285714285714

# --- Selections / active sheet ---
# MT_10's selection grows to include the two neighboring cells.
$ws1.Range("N56:P56").Select()

# Sheet2 becomes the active sheet, with H8:H9 selected (active cell H8).
$ws2.Activate()
$ws2.Range("H8:H9").Select()
